$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the row above (row 16) into the new row 17, so that
# styles (date/time number formats, centered Sno, wrapped description, etc.)
# line up with the existing cellXfs instead of creating new ones.
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the data for the new log entry (row 17)
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 44728
$ws.Cells.Item(17, 3).Value = 0.79166666666666663
$ws.Cells.Item(17, 4).Value = 0.875
$ws.Cells.Item(17, 5).Formula = "=D17-C17"
$ws.Cells.Item(17, 6).Value = "Code"
$ws.Cells.Item(17, 7).Value = "1. resizing BDD100k dataset to 256  x 512 size image`n2. file was very large,so had to use kaggle to download zip files and converting to png images`n3. labels conversion pending,yet to be done"

# Row 17 holds a 3-line wrapped description, so it needs a taller row (like row 16)
$ws.Rows.Item(17).RowHeight = 60

# Recalculate so the SUM total in row 22 (E22) picks up the new row's hours
$excel.Calculate()

# Move/save the active selection like the author did after logging the entry
$ws.Range("G18").Select()
